{"js": "// The document contains a single 20x5 table of simple arithmetic\n// problems (e.g. \"84-33=\"). The commit replaces the text of 100 of\n// these cells (in row-major reading order) with new problem strings,\n// leaving every other part of the document (run formatting, paragraph\n// formatting, table structure, etc.) untouched.\nconst newValues = [\"29+33=\",\"72-58=\",\"26-7=\",\"20+73=\",\"20+51=\",\"41+31=\",\"38+27=\",\"63-63=\",\"49-45=\",\"18+72=\",\"48-12=\",\"34-21=\",\"69+6=\",\"24+26=\",\"61-42=\",\"4+78=\",\"10+29=\",\"66+0=\",\"16-3=\",\"9+58=\",\"42-16=\",\"14+1=\",\"10+78=\",\"99-92=\",\"35-28=\",\"16+41=\",\"95-19=\",\"22+73=\",\"99-41=\",\"5+62=\",\"75-71=\",\"62+27=\",\"67-13=\",\"1+62=\",\"57-50=\",\"18+65=\",\"0+96=\",\"95-68=\",\"39-33=\",\"53+16=\",\"19-6=\",\"85-53=\",\"4+78=\",\"0+49=\",\"26+46=\",\"86-58=\",\"22-13=\",\"62-41=\",\"48-42=\",\"72-48=\",\"38-19=\",\"13-4=\",\"93-74=\",\"98-88=\",\"4+70=\",\"75+4=\",\"62-10=\",\"86-50=\",\"13+78=\",\"4+86=\",\"83-65=\",\"26+43=\",\"76-7=\",\"81-47=\",\"35-19=\",\"40+17=\",\"26-24=\",\"85-60=\",\"89-27=\",\"45-14=\",\"14+25=\",\"49+9=\",\"52-1=\",\"28+38=\",\"91-43=\",\"9+80=\",\"89-26=\",\"84-16=\",\"44+42=\",\"18+39=\",\"43-24=\",\"75-49=\",\"40-28=\",\"47-42=\",\"38-21=\",\"22+15=\",\"2+53=\",\"51-8=\",\"14+31=\",\"17+48=\",\"31+66=\",\"56+0=\",\"98-39=\",\"13+3=\",\"90-34=\",\"59+4=\",\"29-3=\",\"4+82=\",\"72+7=\",\"22-17=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst colCount = 5;\nconst rowCount = table.rowCount;\nif (rowCount * colCount !== newValues.length) {\n  throw new Error(`Unexpected table shape ${rowCount}x${colCount}, expected ${newValues.length} cells`);\n}\n\n// Collect the first-paragraph range of every cell so we can overwrite\n// just the run text (InsertLocation.Replace on an existing range keeps\n// that run's formatting) instead of replacing the whole cell body\n// (which would reset rPr/pPr to cell defaults).\nconst ranges = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.paragraphs.getFirst().getRange();\n    ranges.push(range);\n  }\n}\n\nfor (let i = 0; i < ranges.length; i++) {\n  ranges[i].insertText(newValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$newValues = @(\n    \"29+33=\",\n    \"72-58=\",\n    \"26-7=\",\n    \"20+73=\",\n    \"20+51=\",\n    \"41+31=\",\n    \"38+27=\",\n    \"63-63=\",\n    \"49-45=\",\n    \"18+72=\",\n    \"48-12=\",\n    \"34-21=\",\n    \"69+6=\",\n    \"24+26=\",\n    \"61-42=\",\n    \"4+78=\",\n    \"10+29=\",\n    \"66+0=\",\n    \"16-3=\",\n    \"9+58=\",\n    \"42-16=\",\n    \"14+1=\",\n    \"10+78=\",\n    \"99-92=\",\n    \"35-28=\",\n    \"16+41=\",\n    \"95-19=\",\n    \"22+73=\",\n    \"99-41=\",\n    \"5+62=\",\n    \"75-71=\",\n    \"62+27=\",\n    \"67-13=\",\n    \"1+62=\",\n    \"57-50=\",\n    \"18+65=\",\n    \"0+96=\",\n    \"95-68=\",\n    \"39-33=\",\n    \"53+16=\",\n    \"19-6=\",\n    \"85-53=\",\n    \"4+78=\",\n    \"0+49=\",\n    \"26+46=\",\n    \"86-58=\",\n    \"22-13=\",\n    \"62-41=\",\n    \"48-42=\",\n    \"72-48=\",\n    \"38-19=\",\n    \"13-4=\",\n    \"93-74=\",\n    \"98-88=\",\n    \"4+70=\",\n    \"75+4=\",\n    \"62-10=\",\n    \"86-50=\",\n    \"13+78=\",\n    \"4+86=\",\n    \"83-65=\",\n    \"26+43=\",\n    \"76-7=\",\n    \"81-47=\",\n    \"35-19=\",\n    \"40+17=\",\n    \"26-24=\",\n    \"85-60=\",\n    \"89-27=\",\n    \"45-14=\",\n    \"14+25=\",\n    \"49+9=\",\n    \"52-1=\",\n    \"28+38=\",\n    \"91-43=\",\n    \"9+80=\",\n    \"89-26=\",\n    \"84-16=\",\n    \"44+42=\",\n    \"18+39=\",\n    \"43-24=\",\n    \"75-49=\",\n    \"40-28=\",\n    \"47-42=\",\n    \"38-21=\",\n    \"22+15=\",\n    \"2+53=\",\n    \"51-8=\",\n    \"14+31=\",\n    \"17+48=\",\n    \"31+66=\",\n    \"56+0=\",\n    \"98-39=\",\n    \"13+3=\",\n    \"90-34=\",\n    \"59+4=\",\n    \"29-3=\",\n    \"4+82=\",\n    \"72+7=\",\n    \"22-17=\"\n)\n\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\nif (($rowCount * $colCount) -ne $newValues.Count) {\n    throw \"Unexpected table shape $rowCount x $colCount, expected $($newValues.Count) cells\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
